$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.901.32"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.889.58"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'0.7673"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'242.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'25.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "'0.07141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").Value = "'0.08529"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("D12").Value = "'0.7638"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.899.62"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "'5.368"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Value = "'93.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "'6.158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "29.846.26"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'13.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "'244.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "'0.000007811"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'8.026"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'0.1622"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "'9.406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'163.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "'18.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'1.520"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'4.498"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "'4.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").Value = "'0.05444"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "'0.7454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").Value = "'2.694"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").Value = "'0.01952"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").Value = "'2.780"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.4470"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "1.102.47"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Value = "'73.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "'6.087"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "'0.8523"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'102.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "'1.869"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "'7.668"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").Value = "'3.075"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.002.82"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06088"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
